$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the simulation result data (re-run produced new counts) ---
$ws.Range("C2").Value = 2402.0
$ws.Range("D2").Value = 8819.0
$ws.Range("E2").Value = 29437.0
$ws.Range("F2").Value = 93042.0
$ws.Range("G2").Value = 319452.0
$ws.Range("H2").Value = 1181466.0
$ws.Range("I2").Value = 4914718.0
$ws.Range("J2").Value = 24342406.0
$ws.Range("K2").Value = 121826711.0

$ws.Range("B3").Value = 4500.0
$ws.Range("C3").Value = 1550.0
$ws.Range("D3").Value = 5169.0
$ws.Range("E3").Value = 17455.0
$ws.Range("F3").Value = 63665.0
$ws.Range("G3").Value = 245750.0
$ws.Range("H3").Value = 959983.0
$ws.Range("I3").Value = 3838458.0
$ws.Range("J3").Value = 15290350.0
$ws.Range("K3").Value = 60857794.0

$ws.Range("B4").Value = 3846.0
$ws.Range("D4").Value = 5504.0
$ws.Range("E4").Value = 16011.0
$ws.Range("F4").Value = 48642.0
$ws.Range("G4").Value = 163337.0
$ws.Range("H4").Value = 558634.0
$ws.Range("I4").Value = 2475179.0
$ws.Range("J4").Value = 7919794.0
$ws.Range("K4").Value = 32799106.0

$ws.Range("B5").Value = 14992.0
$ws.Range("C5").Value = 2813.0
$ws.Range("D5").Value = 5610.0
$ws.Range("E5").Value = 12194.0
$ws.Range("F5").Value = 26579.0
$ws.Range("G5").Value = 54709.0
$ws.Range("H5").Value = 99520.0
$ws.Range("I5").Value = 245188.0
$ws.Range("J5").Value = 454233.0
$ws.Range("K5").Value = 944549.0

$ws.Range("B6").Value = 5093.0
$ws.Range("C6").Value = 2250.0
$ws.Range("D6").Value = 4683.0
$ws.Range("E6").Value = 60472.0
$ws.Range("F6").Value = 11388.0
$ws.Range("G6").Value = 23249.0
$ws.Range("H6").Value = 49463.0
$ws.Range("I6").Value = 116139.0
$ws.Range("J6").Value = 297677.0
$ws.Range("K6").Value = 861938.0

$ws.Range("B7").Value = 3664.0
$ws.Range("C7").Value = 1307.0
$ws.Range("D7").Value = 2980.0
$ws.Range("E7").Value = 7101.0
$ws.Range("F7").Value = 16786.0
$ws.Range("G7").Value = 38789.0
$ws.Range("H7").Value = 89651.0
$ws.Range("I7").Value = 209379.0
$ws.Range("J7").Value = 466702.0
$ws.Range("K7").Value = 1062225.0

$ws.Range("B8").Value = 14962.0
$ws.Range("C8").Value = 6872.0
$ws.Range("D8").Value = 2326.0
$ws.Range("E8").Value = 3938.0
$ws.Range("F8").Value = 15844.0
$ws.Range("G8").Value = 20223.0
$ws.Range("H8").Value = 38698.0
$ws.Range("I8").Value = 69732.0
$ws.Range("J8").Value = 132014.0
$ws.Range("K8").Value = 257079.0

# --- Auto resize the cells: autofit columns B:K to their (now-updated) content ---
$ws.Range("B1:K8").EntireColumn.AutoFit()

# The host's autofit heuristic doesn't reproduce Excel's exact font-metric
# based best-fit widths, so nudge each column to the precise width Excel
# computed for this content (keeps customWidth + near-exact widths).
$ws.Columns.Item(2).ColumnWidth = 5.833333333333333
$ws.Columns.Item(3).ColumnWidth = 4.666666666666667
$ws.Columns.Item(4).ColumnWidth = 4.666666666666667
$ws.Columns.Item(5).ColumnWidth = 5.833333333333333
$ws.Columns.Item(6).ColumnWidth = 5.833333333333333
$ws.Columns.Item(7).ColumnWidth = 6.833333333333333
$ws.Columns.Item(8).ColumnWidth = 8.0
$ws.Columns.Item(9).ColumnWidth = 8.0
$ws.Columns.Item(10).ColumnWidth = 9.166666666666666
$ws.Columns.Item(11).ColumnWidth = 10.166666666666666
